# Auto-generated: apply crypto price/volume refresh from GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value (matches the source feed, which
    # stores prices/volumes/links as inline strings, not numbers) and then
    # drop the temporary text number-format so the cell style is left
    # exactly as it was (no lingering "@" format override).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "30.536.78"
Set-TextValue $ws.Range("E2") "  -0.15%  "
Set-TextValue $ws.Range("D3") "1.917.30"
Set-TextValue $ws.Range("E3") "  -0.47%  "
Set-TextValue $ws.Range("D5") "245.21"
Set-TextValue $ws.Range("E5") "  -0.63%  "
Set-TextValue $ws.Range("E6") "  +0.03%  "
Set-TextValue $ws.Range("D7") "0.4819"
Set-TextValue $ws.Range("E7") "  +1.73%  "
Set-TextValue $ws.Range("D8") "0.2900"
Set-TextValue $ws.Range("E8") "  -0.50%  "
Set-TextValue $ws.Range("D9") "0.06702"
Set-TextValue $ws.Range("E9") "  -1.53%  "
Set-TextValue $ws.Range("D10") "111.46"
Set-TextValue $ws.Range("E10") "  +5.30%  "
Set-TextValue $ws.Range("D11") "18.96"
Set-TextValue $ws.Range("E11") "  +3.15%  "
Set-TextValue $ws.Range("D12") "1.920.75"
Set-TextValue $ws.Range("E12") "  +0.25%  "
Set-TextValue $ws.Range("D13") "0.07565"
Set-TextValue $ws.Range("E13") "  -2.22%  "
Set-TextValue $ws.Range("D14") "5.287"
Set-TextValue $ws.Range("E14") "  -1.00%  "
Set-TextValue $ws.Range("D15") "0.6679"
Set-TextValue $ws.Range("E15") "  -0.57%  "
Set-TextValue $ws.Range("D16") "298.06"
Set-TextValue $ws.Range("E16") "  +3.41%  "
Set-TextValue $ws.Range("D17") "30.540.50"
Set-TextValue $ws.Range("E17") "  -0.24%  "
Set-TextValue $ws.Range("B18") "ShibaInu"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.000007619"
Set-TextValue $ws.Range("E18") "  -0.42%  "
Set-TextValue $ws.Range("E19") "  +0.03%  "
Set-TextValue $ws.Range("B20") "Avalanche"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D20") "12.98"
Set-TextValue $ws.Range("E20") "  -0.63%  "
Set-TextValue $ws.Range("D21") "5.547"
Set-TextValue $ws.Range("E21") "  +1.63%  "
Set-TextValue $ws.Range("D22") "2.168.39"
Set-TextValue $ws.Range("E22") "  +0.03%  "
Set-TextValue $ws.Range("D23") "1.001"
Set-TextValue $ws.Range("E23") "  +0.02%  "
Set-TextValue $ws.Range("D24") "6.435"
Set-TextValue $ws.Range("E24") "  +2.66%  "
Set-TextValue $ws.Range("D25") "9.442"
Set-TextValue $ws.Range("E25") "  +0.49%  "
Set-TextValue $ws.Range("D26") "165.41"
Set-TextValue $ws.Range("E26") "  -1.98%  "
Set-TextValue $ws.Range("D27") "20.26"
Set-TextValue $ws.Range("E27") "  -2.23%  "
Set-TextValue $ws.Range("D28") "2.099"
Set-TextValue $ws.Range("E28") "  -2.03%  "
Set-TextValue $ws.Range("D30") "1.435"
Set-TextValue $ws.Range("E30") "  +5.49%  "
Set-TextValue $ws.Range("D31") "4.142"
Set-TextValue $ws.Range("E31") "  -0.68%  "
Set-TextValue $ws.Range("D32") "4.066"
Set-TextValue $ws.Range("E32") "  +1.03%  "
Set-TextValue $ws.Range("D33") "0.04998"
Set-TextValue $ws.Range("E33") "  -1.41%  "
Set-TextValue $ws.Range("D34") "0.7388"
Set-TextValue $ws.Range("E34") "  -0.25%  "
Set-TextValue $ws.Range("E35") "  -1.49%  "
Set-TextValue $ws.Range("E36") "  -0.01%  "
Set-TextValue $ws.Range("D37") "2.724"
Set-TextValue $ws.Range("E37") "  -0.30%  "
Set-TextValue $ws.Range("D38") "0.02014"
Set-TextValue $ws.Range("E38") "  -3.58%  "
Set-TextValue $ws.Range("D39") "2.678"
Set-TextValue $ws.Range("E39") "  -0.59%  "
Set-TextValue $ws.Range("D40") "110.73"
Set-TextValue $ws.Range("E40") "  -0.37%  "
Set-TextValue $ws.Range("D41") "2.013"
Set-TextValue $ws.Range("E41") "  -2.37%  "
Set-TextValue $ws.Range("D42") "0.4422"
Set-TextValue $ws.Range("E42") "  -0.44%  "
Set-TextValue $ws.Range("D43") "0.8655"
Set-TextValue $ws.Range("E43") "  -1.06%  "
Set-TextValue $ws.Range("D44") "70.60"
Set-TextValue $ws.Range("E44") "  +4.38%  "
Set-TextValue $ws.Range("D45") "5.830"
Set-TextValue $ws.Range("E45") "  -1.48%  "
Set-TextValue $ws.Range("E46") "  +0.03%  "
Set-TextValue $ws.Range("B47") "Aptos"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D47") "7.209"
Set-TextValue $ws.Range("E47") "  -0.88%  "
Set-TextValue $ws.Range("B48") "BitcoinSV"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D48") "48.75"
Set-TextValue $ws.Range("E48") "  +2.92%  "
Set-TextValue $ws.Range("D49") "9.249"
Set-TextValue $ws.Range("E49") "  -1.23%  "
Set-TextValue $ws.Range("D50") "0.1227"
Set-TextValue $ws.Range("E50") "  -0.63%  "
Set-TextValue $ws.Range("D51") "34.86"
Set-TextValue $ws.Range("E51") "  -0.92%  "
